$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 91, pushing rows 91..145 down to 92..146.
# Excel copies the formatting of the row above (row 90) to the new row,
# including the date number-format on column D.
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new "Ralli Seedless" record.
$ws.Cells.Item(91, 1).Value = 8
$ws.Cells.Item(91, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(91, 3).Value = "Coquimbo"
$ws.Cells.Item(91, 4).Value = 44981
$ws.Cells.Item(91, 5).Value = 4
$ws.Cells.Item(91, 6).Value = "Fruta"
$ws.Cells.Item(91, 7).Value = 100109
$ws.Cells.Item(91, 8).Value = "Uva"
$ws.Cells.Item(91, 9).Value = 100109001
$ws.Cells.Item(91, 10).Value = "Uva"
$ws.Cells.Item(91, 11).Value = "Ralli Seedless"
$ws.Cells.Item(91, 12).Value = "Primera"
$ws.Cells.Item(91, 13).Value = 300
$ws.Cells.Item(91, 14).Value = 10000
$ws.Cells.Item(91, 15).Value = 11000
$ws.Cells.Item(91, 16).Value = 10500
$ws.Cells.Item(91, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(91, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(91, 19).Value = 583
$ws.Cells.Item(91, 20).Value = 18
